$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: 'A 9538-2024' (source row 2)
$ws.Range("A2").Value = 'A 9538-2024'
$ws.Range("B2").Value = 45359
$ws.Range("C2").Value = 46064
$ws.Range("D2").Value = 'STOCKHOLMS LÄN'
$ws.Range("E2").Value = 'SALEM'
$ws.Range("F2").Value = 'Kommuner'
$ws.Range("G2").Value = 14
$ws.Range("H2").Value = 6
$ws.Range("I2").Value = 15
$ws.Range("J2").Value = 12
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 14
$ws.Range("P2").Value = 2
$ws.Range("Q2").Value = 32
$ws.Range("R2").Value = 'Läderdoftande fingersvamp
Vågticka
Duvhök
Fyrflikig jordstjärna
Grön aspvedbock
Gul taggsvamp
Gultoppig fingersvamp
Motaggsvamp
Reliktbock
Spillkråka
Tallticka
Talltita
Ullticka
Vedskivlav
Björksplintborre
Blomkålssvamp
Blåmossa
Brandticka
Bronshjon
Dropptaggsvamp
Granbarkgnagare
Grovticka
Grönpyrola
Jättesvampmal
Mindre märgborre
Mörk husmossa
Skarp dropptaggsvamp
Vedticka
Vågbandad barkbock
Grönsiska
Kungsfågel
Revlummer'
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/artfynd/A 9538-2024 artfynd.xlsx", "A 9538-2024")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/kartor/A 9538-2024 karta.png", "A 9538-2024")'
$ws.Range("U2").ClearContents()
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/klagomål/A 9538-2024 FSC-klagomål.docx", "A 9538-2024")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/klagomålsmail/A 9538-2024 FSC-klagomål mail.docx", "A 9538-2024")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/tillsyn/A 9538-2024 tillsynsbegäran.docx", "A 9538-2024")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/tillsynsmail/A 9538-2024 tillsynsbegäran mail.docx", "A 9538-2024")'
$ws.Range("Z2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/fåglar/A 9538-2024 prioriterade fågelarter.docx", "A 9538-2024")'

# Row 3: 'A 6162-2025' (source row 3)
$ws.Range("A3").Value = 'A 6162-2025'
$ws.Range("B3").Value = 45698
$ws.Range("C3").Value = 46064
$ws.Range("D3").Value = 'STOCKHOLMS LÄN'
$ws.Range("E3").Value = 'SALEM'
$ws.Range("F3").ClearContents()
$ws.Range("G3").Value = 22.3
$ws.Range("H3").Value = 6
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 4
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 5
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 8
$ws.Range("R3").Value = 'Knärot
Entita
Spillkråka
Tallticka
Talltita
Blåmossa
Kungsfågel
Tjäder'
$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/artfynd/A 6162-2025 artfynd.xlsx", "A 6162-2025")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/kartor/A 6162-2025 karta.png", "A 6162-2025")'
$ws.Range("U3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/knärot/A 6162-2025 karta knärot.png", "A 6162-2025")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/klagomål/A 6162-2025 FSC-klagomål.docx", "A 6162-2025")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/klagomålsmail/A 6162-2025 FSC-klagomål mail.docx", "A 6162-2025")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/tillsyn/A 6162-2025 tillsynsbegäran.docx", "A 6162-2025")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/tillsynsmail/A 6162-2025 tillsynsbegäran mail.docx", "A 6162-2025")'
$ws.Range("Z3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/fåglar/A 6162-2025 prioriterade fågelarter.docx", "A 6162-2025")'

# Row 4: 'A 36244-2024' (source row 4)
$ws.Range("A4").Value = 'A 36244-2024'
$ws.Range("B4").Value = 45534
$ws.Range("C4").Value = 46064
$ws.Range("D4").Value = 'STOCKHOLMS LÄN'
$ws.Range("E4").Value = 'SALEM'
$ws.Range("F4").Value = 'Kommuner'
$ws.Range("G4").Value = 0.1
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 4
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 2
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 6
$ws.Range("R4").Value = 'Kandelabersvamp
Ullticka
Fällmossa
Grov fjädermossa
Platt fjädermossa
Västlig hakmossa'
$ws.Range("S4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/artfynd/A 36244-2024 artfynd.xlsx", "A 36244-2024")'
$ws.Range("T4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/kartor/A 36244-2024 karta.png", "A 36244-2024")'
$ws.Range("U4").ClearContents()
$ws.Range("V4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/klagomål/A 36244-2024 FSC-klagomål.docx", "A 36244-2024")'
$ws.Range("W4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/klagomålsmail/A 36244-2024 FSC-klagomål mail.docx", "A 36244-2024")'
$ws.Range("X4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/tillsyn/A 36244-2024 tillsynsbegäran.docx", "A 36244-2024")'
$ws.Range("Y4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/tillsynsmail/A 36244-2024 tillsynsbegäran mail.docx", "A 36244-2024")'
$ws.Range("Z4").ClearContents()

# Row 5: 'A 61302-2022' (source row 6)
$ws.Range("A5").Value = 'A 61302-2022'
$ws.Range("B5").Value = 44915
$ws.Range("C5").Value = 46064
$ws.Range("D5").Value = 'STOCKHOLMS LÄN'
$ws.Range("E5").Value = 'SALEM'
$ws.Range("F5").Value = 'Kommuner'
$ws.Range("G5").Value = 14.8
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 2
$ws.Range("R5").Value = 'Trådticka
Zontaggsvamp'
$ws.Range("S5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/artfynd/A 61302-2022 artfynd.xlsx", "A 61302-2022")'
$ws.Range("T5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/kartor/A 61302-2022 karta.png", "A 61302-2022")'
$ws.Range("U5").ClearContents()
$ws.Range("V5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/klagomål/A 61302-2022 FSC-klagomål.docx", "A 61302-2022")'
$ws.Range("W5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/klagomålsmail/A 61302-2022 FSC-klagomål mail.docx", "A 61302-2022")'
$ws.Range("X5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/tillsyn/A 61302-2022 tillsynsbegäran.docx", "A 61302-2022")'
$ws.Range("Y5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/tillsynsmail/A 61302-2022 tillsynsbegäran mail.docx", "A 61302-2022")'
$ws.Range("Z5").ClearContents()

# Row 6: 'A 2992-2023' (source row 5)
$ws.Range("A6").Value = 'A 2992-2023'
$ws.Range("B6").Value = 44945
$ws.Range("C6").Value = 46064
$ws.Range("D6").Value = 'STOCKHOLMS LÄN'
$ws.Range("E6").Value = 'SALEM'
$ws.Range("F6").Value = 'Kommuner'
$ws.Range("G6").Value = 1.7
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 1
$ws.Range("P6").Value = 1
$ws.Range("Q6").Value = 2
$ws.Range("R6").Value = 'Rynkskinn
Blåmossa'
$ws.Range("S6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/artfynd/A 2992-2023 artfynd.xlsx", "A 2992-2023")'
$ws.Range("T6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/kartor/A 2992-2023 karta.png", "A 2992-2023")'
$ws.Range("U6").ClearContents()
$ws.Range("V6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/klagomål/A 2992-2023 FSC-klagomål.docx", "A 2992-2023")'
$ws.Range("W6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/klagomålsmail/A 2992-2023 FSC-klagomål mail.docx", "A 2992-2023")'
$ws.Range("X6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/tillsyn/A 2992-2023 tillsynsbegäran.docx", "A 2992-2023")'
$ws.Range("Y6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/tillsynsmail/A 2992-2023 tillsynsbegäran mail.docx", "A 2992-2023")'
$ws.Range("Z6").ClearContents()

# Row 7: 'A 68700-2021' (source row 7)
$ws.Range("A7").Value = 'A 68700-2021'
$ws.Range("B7").Value = 44529
$ws.Range("C7").Value = 46064
$ws.Range("D7").Value = 'STOCKHOLMS LÄN'
$ws.Range("E7").Value = 'SALEM'
$ws.Range("F7").Value = 'Kommuner'
$ws.Range("G7").Value = 5.3
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = 'Blåmossa'
$ws.Range("S7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/artfynd/A 68700-2021 artfynd.xlsx", "A 68700-2021")'
$ws.Range("T7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/kartor/A 68700-2021 karta.png", "A 68700-2021")'
$ws.Range("U7").ClearContents()
$ws.Range("V7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/klagomål/A 68700-2021 FSC-klagomål.docx", "A 68700-2021")'
$ws.Range("W7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/klagomålsmail/A 68700-2021 FSC-klagomål mail.docx", "A 68700-2021")'
$ws.Range("X7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/tillsyn/A 68700-2021 tillsynsbegäran.docx", "A 68700-2021")'
$ws.Range("Y7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/tillsynsmail/A 68700-2021 tillsynsbegäran mail.docx", "A 68700-2021")'
$ws.Range("Z7").ClearContents()

# Row 8: 'A 47262-2024' (source row 10)
$ws.Range("A8").Value = 'A 47262-2024'
$ws.Range("B8").Value = 45586
$ws.Range("C8").Value = 46064
$ws.Range("D8").Value = 'STOCKHOLMS LÄN'
$ws.Range("E8").Value = 'SALEM'
$ws.Range("F8").ClearContents()
$ws.Range("G8").Value = 2
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 1
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = 1
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 1
$ws.Range("R8").Value = 'Spillkråka'
$ws.Range("S8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/artfynd/A 47262-2024 artfynd.xlsx", "A 47262-2024")'
$ws.Range("T8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/kartor/A 47262-2024 karta.png", "A 47262-2024")'
$ws.Range("U8").ClearContents()
$ws.Range("V8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/klagomål/A 47262-2024 FSC-klagomål.docx", "A 47262-2024")'
$ws.Range("W8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/klagomålsmail/A 47262-2024 FSC-klagomål mail.docx", "A 47262-2024")'
$ws.Range("X8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/tillsyn/A 47262-2024 tillsynsbegäran.docx", "A 47262-2024")'
$ws.Range("Y8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/tillsynsmail/A 47262-2024 tillsynsbegäran mail.docx", "A 47262-2024")'
$ws.Range("Z8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/fåglar/A 47262-2024 prioriterade fågelarter.docx", "A 47262-2024")'

# Row 9: 'A 33411-2025' (source row 8)
$ws.Range("A9").Value = 'A 33411-2025'
$ws.Range("B9").Value = 45839
$ws.Range("C9").Value = 46064
$ws.Range("D9").Value = 'STOCKHOLMS LÄN'
$ws.Range("E9").Value = 'SALEM'
$ws.Range("F9").ClearContents()
$ws.Range("G9").Value = 2.8
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 1
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0
$ws.Range("O9").Value = 1
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = 'Spillkråka'
$ws.Range("S9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/artfynd/A 33411-2025 artfynd.xlsx", "A 33411-2025")'
$ws.Range("T9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/kartor/A 33411-2025 karta.png", "A 33411-2025")'
$ws.Range("U9").ClearContents()
$ws.Range("V9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/klagomål/A 33411-2025 FSC-klagomål.docx", "A 33411-2025")'
$ws.Range("W9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/klagomålsmail/A 33411-2025 FSC-klagomål mail.docx", "A 33411-2025")'
$ws.Range("X9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/tillsyn/A 33411-2025 tillsynsbegäran.docx", "A 33411-2025")'
$ws.Range("Y9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/tillsynsmail/A 33411-2025 tillsynsbegäran mail.docx", "A 33411-2025")'
$ws.Range("Z9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/fåglar/A 33411-2025 prioriterade fågelarter.docx", "A 33411-2025")'

# Row 10: 'A 32951-2025' (source row 9)
$ws.Range("A10").Value = 'A 32951-2025'
$ws.Range("B10").Value = 45839
$ws.Range("C10").Value = 46064
$ws.Range("D10").Value = 'STOCKHOLMS LÄN'
$ws.Range("E10").Value = 'SALEM'
$ws.Range("F10").ClearContents()
$ws.Range("G10").Value = 2.6
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 1
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 1
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 1
$ws.Range("R10").Value = 'Etternässla'
$ws.Range("S10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/artfynd/A 32951-2025 artfynd.xlsx", "A 32951-2025")'
$ws.Range("T10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/kartor/A 32951-2025 karta.png", "A 32951-2025")'
$ws.Range("U10").ClearContents()
$ws.Range("V10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/klagomål/A 32951-2025 FSC-klagomål.docx", "A 32951-2025")'
$ws.Range("W10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/klagomålsmail/A 32951-2025 FSC-klagomål mail.docx", "A 32951-2025")'
$ws.Range("X10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/tillsyn/A 32951-2025 tillsynsbegäran.docx", "A 32951-2025")'
$ws.Range("Y10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/tillsynsmail/A 32951-2025 tillsynsbegäran mail.docx", "A 32951-2025")'
$ws.Range("Z10").ClearContents()

# Row 11: 'A 2987-2023' (source row 11)
$ws.Range("A11").Value = 'A 2987-2023'
$ws.Range("B11").Value = 44945
$ws.Range("C11").Value = 46064
$ws.Range("D11").Value = 'STOCKHOLMS LÄN'
$ws.Range("E11").Value = 'SALEM'
$ws.Range("F11").Value = 'Kommuner'
$ws.Range("G11").Value = 6.2
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 1
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 1
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 1
$ws.Range("R11").Value = 'Tallticka'
$ws.Range("S11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/artfynd/A 2987-2023 artfynd.xlsx", "A 2987-2023")'
$ws.Range("T11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/kartor/A 2987-2023 karta.png", "A 2987-2023")'
$ws.Range("U11").ClearContents()
$ws.Range("V11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/klagomål/A 2987-2023 FSC-klagomål.docx", "A 2987-2023")'
$ws.Range("W11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/klagomålsmail/A 2987-2023 FSC-klagomål mail.docx", "A 2987-2023")'
$ws.Range("X11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/tillsyn/A 2987-2023 tillsynsbegäran.docx", "A 2987-2023")'
$ws.Range("Y11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0128/tillsynsmail/A 2987-2023 tillsynsbegäran mail.docx", "A 2987-2023")'
$ws.Range("Z11").ClearContents()

# Row 12: 'A 28838-2021' (source row 12)
$ws.Range("A12").Value = 'A 28838-2021'
$ws.Range("B12").Value = 44357
$ws.Range("C12").Value = 46064
$ws.Range("D12").Value = 'STOCKHOLMS LÄN'
$ws.Range("E12").Value = 'SALEM'
$ws.Range("F12").ClearContents()
$ws.Range("G12").Value = 6.4
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 0
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 0
$ws.Range("R12").Value = ''
$ws.Range("S12").ClearContents()
$ws.Range("T12").ClearContents()
$ws.Range("U12").ClearContents()
$ws.Range("V12").ClearContents()
$ws.Range("W12").ClearContents()
$ws.Range("X12").ClearContents()
$ws.Range("Y12").ClearContents()
$ws.Range("Z12").ClearContents()

# Row 13: 'A 13351-2021' (source row 13)
$ws.Range("A13").Value = 'A 13351-2021'
$ws.Range("B13").Value = 44273
$ws.Range("C13").Value = 46064
$ws.Range("D13").Value = 'STOCKHOLMS LÄN'
$ws.Range("E13").Value = 'SALEM'
$ws.Range("F13").ClearContents()
$ws.Range("G13").Value = 2.6
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 0
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 0
$ws.Range("R13").Value = ''
$ws.Range("S13").ClearContents()
$ws.Range("T13").ClearContents()
$ws.Range("U13").ClearContents()
$ws.Range("V13").ClearContents()
$ws.Range("W13").ClearContents()
$ws.Range("X13").ClearContents()
$ws.Range("Y13").ClearContents()
$ws.Range("Z13").ClearContents()

# Row 14: 'A 43893-2021' (source row 14)
$ws.Range("A14").Value = 'A 43893-2021'
$ws.Range("B14").Value = 44434
$ws.Range("C14").Value = 46064
$ws.Range("D14").Value = 'STOCKHOLMS LÄN'
$ws.Range("E14").Value = 'SALEM'
$ws.Range("F14").Value = 'Kommuner'
$ws.Range("G14").Value = 2.2
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 0
$ws.Range("R14").Value = ''
$ws.Range("S14").ClearContents()
$ws.Range("T14").ClearContents()
$ws.Range("U14").ClearContents()
$ws.Range("V14").ClearContents()
$ws.Range("W14").ClearContents()
$ws.Range("X14").ClearContents()
$ws.Range("Y14").ClearContents()
$ws.Range("Z14").ClearContents()

# Row 15: 'A 51720-2023' (source row 48)
$ws.Range("A15").Value = 'A 51720-2023'
$ws.Range("B15").Value = 45222
$ws.Range("C15").Value = 46064
$ws.Range("D15").Value = 'STOCKHOLMS LÄN'
$ws.Range("E15").Value = 'SALEM'
$ws.Range("F15").ClearContents()
$ws.Range("G15").Value = 0.4
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 0
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = 0
$ws.Range("Q15").Value = 0
$ws.Range("R15").Value = ''
$ws.Range("S15").ClearContents()
$ws.Range("T15").ClearContents()
$ws.Range("U15").ClearContents()
$ws.Range("V15").ClearContents()
$ws.Range("W15").ClearContents()
$ws.Range("X15").ClearContents()
$ws.Range("Y15").ClearContents()
$ws.Range("Z15").ClearContents()

# Row 16: 'A 33262-2021' (source row 45)
$ws.Range("A16").Value = 'A 33262-2021'
$ws.Range("B16").Value = 44377.35033564815
$ws.Range("C16").Value = 46064
$ws.Range("D16").Value = 'STOCKHOLMS LÄN'
$ws.Range("E16").Value = 'SALEM'
$ws.Range("F16").ClearContents()
$ws.Range("G16").Value = 0.7
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = 0
$ws.Range("Q16").Value = 0
$ws.Range("R16").Value = ''
$ws.Range("S16").ClearContents()
$ws.Range("T16").ClearContents()
$ws.Range("U16").ClearContents()
$ws.Range("V16").ClearContents()
$ws.Range("W16").ClearContents()
$ws.Range("X16").ClearContents()
$ws.Range("Y16").ClearContents()
$ws.Range("Z16").ClearContents()

# Row 17: 'A 8300-2024' (source row 16)
$ws.Range("A17").Value = 'A 8300-2024'
$ws.Range("B17").Value = 45352.44717592592
$ws.Range("C17").Value = 46064
$ws.Range("D17").Value = 'STOCKHOLMS LÄN'
$ws.Range("E17").Value = 'SALEM'
$ws.Range("F17").Value = 'Kommuner'
$ws.Range("G17").Value = 3.4
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = 0
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = 0
$ws.Range("Q17").Value = 0
$ws.Range("R17").Value = ''
$ws.Range("S17").ClearContents()
$ws.Range("T17").ClearContents()
$ws.Range("U17").ClearContents()
$ws.Range("V17").ClearContents()
$ws.Range("W17").ClearContents()
$ws.Range("X17").ClearContents()
$ws.Range("Y17").ClearContents()
$ws.Range("Z17").ClearContents()

# Row 18: 'A 1071-2023' (source row 56)
$ws.Range("A18").Value = 'A 1071-2023'
$ws.Range("B18").Value = 44935
$ws.Range("C18").Value = 46064
$ws.Range("D18").Value = 'STOCKHOLMS LÄN'
$ws.Range("E18").Value = 'SALEM'
$ws.Range("F18").Value = 'Kommuner'
$ws.Range("G18").Value = 5.1
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = 0
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = 0
$ws.Range("Q18").Value = 0
$ws.Range("R18").Value = ''
$ws.Range("S18").ClearContents()
$ws.Range("T18").ClearContents()
$ws.Range("U18").ClearContents()
$ws.Range("V18").ClearContents()
$ws.Range("W18").ClearContents()
$ws.Range("X18").ClearContents()
$ws.Range("Y18").ClearContents()
$ws.Range("Z18").ClearContents()

# Row 19: 'A 12874-2021' (source row 17)
$ws.Range("A19").Value = 'A 12874-2021'
$ws.Range("B19").Value = 44271.42787037037
$ws.Range("C19").Value = 46064
$ws.Range("D19").Value = 'STOCKHOLMS LÄN'
$ws.Range("E19").Value = 'SALEM'
$ws.Range("F19").ClearContents()
$ws.Range("G19").Value = 6
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = 0
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = 0
$ws.Range("Q19").Value = 0
$ws.Range("R19").Value = ''
$ws.Range("S19").ClearContents()
$ws.Range("T19").ClearContents()
$ws.Range("U19").ClearContents()
$ws.Range("V19").ClearContents()
$ws.Range("W19").ClearContents()
$ws.Range("X19").ClearContents()
$ws.Range("Y19").ClearContents()
$ws.Range("Z19").ClearContents()

# Row 20: 'A 32962-2025' (source row 35)
$ws.Range("A20").Value = 'A 32962-2025'
$ws.Range("B20").Value = 45839
$ws.Range("C20").Value = 46064
$ws.Range("D20").Value = 'STOCKHOLMS LÄN'
$ws.Range("E20").Value = 'SALEM'
$ws.Range("F20").ClearContents()
$ws.Range("G20").Value = 0.8
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = 0
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = 0
$ws.Range("Q20").Value = 0
$ws.Range("R20").Value = ''
$ws.Range("S20").ClearContents()
$ws.Range("T20").ClearContents()
$ws.Range("U20").ClearContents()
$ws.Range("V20").ClearContents()
$ws.Range("W20").ClearContents()
$ws.Range("X20").ClearContents()
$ws.Range("Y20").ClearContents()
$ws.Range("Z20").ClearContents()

# Row 21: 'A 33029-2025' (source row 28)
$ws.Range("A21").Value = 'A 33029-2025'
$ws.Range("B21").Value = 45839
$ws.Range("C21").Value = 46064
$ws.Range("D21").Value = 'STOCKHOLMS LÄN'
$ws.Range("E21").Value = 'SALEM'
$ws.Range("F21").ClearContents()
$ws.Range("G21").Value = 1.5
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = 0
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = 0
$ws.Range("Q21").Value = 0
$ws.Range("R21").Value = ''
$ws.Range("S21").ClearContents()
$ws.Range("T21").ClearContents()
$ws.Range("U21").ClearContents()
$ws.Range("V21").ClearContents()
$ws.Range("W21").ClearContents()
$ws.Range("X21").ClearContents()
$ws.Range("Y21").ClearContents()
$ws.Range("Z21").ClearContents()

# Row 22: 'A 33455-2025' (source row 23)
$ws.Range("A22").Value = 'A 33455-2025'
$ws.Range("B22").Value = 45839
$ws.Range("C22").Value = 46064
$ws.Range("D22").Value = 'STOCKHOLMS LÄN'
$ws.Range("E22").Value = 'SALEM'
$ws.Range("F22").ClearContents()
$ws.Range("G22").Value = 0.9
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = 0
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = 0
$ws.Range("Q22").Value = 0
$ws.Range("R22").Value = ''
$ws.Range("S22").ClearContents()
$ws.Range("T22").ClearContents()
$ws.Range("U22").ClearContents()
$ws.Range("V22").ClearContents()
$ws.Range("W22").ClearContents()
$ws.Range("X22").ClearContents()
$ws.Range("Y22").ClearContents()
$ws.Range("Z22").ClearContents()

# Row 23: 'A 33243-2025' (source row 24)
$ws.Range("A23").Value = 'A 33243-2025'
$ws.Range("B23").Value = 45839
$ws.Range("C23").Value = 46064
$ws.Range("D23").Value = 'STOCKHOLMS LÄN'
$ws.Range("E23").Value = 'SALEM'
$ws.Range("F23").ClearContents()
$ws.Range("G23").Value = 2.7
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0
$ws.Range("N23").Value = 0
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = 0
$ws.Range("Q23").Value = 0
$ws.Range("R23").Value = ''
$ws.Range("S23").ClearContents()
$ws.Range("T23").ClearContents()
$ws.Range("U23").ClearContents()
$ws.Range("V23").ClearContents()
$ws.Range("W23").ClearContents()
$ws.Range("X23").ClearContents()
$ws.Range("Y23").ClearContents()
$ws.Range("Z23").ClearContents()

# Row 24: 'A 33251-2025' (source row 25)
$ws.Range("A24").Value = 'A 33251-2025'
$ws.Range("B24").Value = 45839
$ws.Range("C24").Value = 46064
$ws.Range("D24").Value = 'STOCKHOLMS LÄN'
$ws.Range("E24").Value = 'SALEM'
$ws.Range("F24").ClearContents()
$ws.Range("G24").Value = 2.9
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 0
$ws.Range("N24").Value = 0
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = 0
$ws.Range("Q24").Value = 0
$ws.Range("R24").Value = ''
$ws.Range("S24").ClearContents()
$ws.Range("T24").ClearContents()
$ws.Range("U24").ClearContents()
$ws.Range("V24").ClearContents()
$ws.Range("W24").ClearContents()
$ws.Range("X24").ClearContents()
$ws.Range("Y24").ClearContents()
$ws.Range("Z24").ClearContents()

# Row 25: 'A 33277-2025' (source row 22)
$ws.Range("A25").Value = 'A 33277-2025'
$ws.Range("B25").Value = 45839
$ws.Range("C25").Value = 46064
$ws.Range("D25").Value = 'STOCKHOLMS LÄN'
$ws.Range("E25").Value = 'SALEM'
$ws.Range("F25").ClearContents()
$ws.Range("G25").Value = 2.4
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("N25").Value = 0
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = 0
$ws.Range("Q25").Value = 0
$ws.Range("R25").Value = ''
$ws.Range("S25").ClearContents()
$ws.Range("T25").ClearContents()
$ws.Range("U25").ClearContents()
$ws.Range("V25").ClearContents()
$ws.Range("W25").ClearContents()
$ws.Range("X25").ClearContents()
$ws.Range("Y25").ClearContents()
$ws.Range("Z25").ClearContents()

# Row 26: 'A 33284-2025' (source row 34)
$ws.Range("A26").Value = 'A 33284-2025'
$ws.Range("B26").Value = 45839
$ws.Range("C26").Value = 46064
$ws.Range("D26").Value = 'STOCKHOLMS LÄN'
$ws.Range("E26").Value = 'SALEM'
$ws.Range("F26").ClearContents()
$ws.Range("G26").Value = 3.3
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 0
$ws.Range("N26").Value = 0
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = 0
$ws.Range("Q26").Value = 0
$ws.Range("R26").Value = ''
$ws.Range("S26").ClearContents()
$ws.Range("T26").ClearContents()
$ws.Range("U26").ClearContents()
$ws.Range("V26").ClearContents()
$ws.Range("W26").ClearContents()
$ws.Range("X26").ClearContents()
$ws.Range("Y26").ClearContents()
$ws.Range("Z26").ClearContents()

# Row 27: 'A 32960-2025' (source row 29)
$ws.Range("A27").Value = 'A 32960-2025'
$ws.Range("B27").Value = 45839
$ws.Range("C27").Value = 46064
$ws.Range("D27").Value = 'STOCKHOLMS LÄN'
$ws.Range("E27").Value = 'SALEM'
$ws.Range("F27").ClearContents()
$ws.Range("G27").Value = 5.3
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = 0
$ws.Range("N27").Value = 0
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = 0
$ws.Range("Q27").Value = 0
$ws.Range("R27").Value = ''
$ws.Range("S27").ClearContents()
$ws.Range("T27").ClearContents()
$ws.Range("U27").ClearContents()
$ws.Range("V27").ClearContents()
$ws.Range("W27").ClearContents()
$ws.Range("X27").ClearContents()
$ws.Range("Y27").ClearContents()
$ws.Range("Z27").ClearContents()

# Row 28: 'A 32968-2025' (source row 30)
$ws.Range("A28").Value = 'A 32968-2025'
$ws.Range("B28").Value = 45839
$ws.Range("C28").Value = 46064
$ws.Range("D28").Value = 'STOCKHOLMS LÄN'
$ws.Range("E28").Value = 'SALEM'
$ws.Range("F28").ClearContents()
$ws.Range("G28").Value = 1.1
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0
$ws.Range("N28").Value = 0
$ws.Range("O28").Value = 0
$ws.Range("P28").Value = 0
$ws.Range("Q28").Value = 0
$ws.Range("R28").Value = ''
$ws.Range("S28").ClearContents()
$ws.Range("T28").ClearContents()
$ws.Range("U28").ClearContents()
$ws.Range("V28").ClearContents()
$ws.Range("W28").ClearContents()
$ws.Range("X28").ClearContents()
$ws.Range("Y28").ClearContents()
$ws.Range("Z28").ClearContents()

# Row 29: 'A 32942-2025' (source row 36)
$ws.Range("A29").Value = 'A 32942-2025'
$ws.Range("B29").Value = 45839
$ws.Range("C29").Value = 46064
$ws.Range("D29").Value = 'STOCKHOLMS LÄN'
$ws.Range("E29").Value = 'SALEM'
$ws.Range("F29").ClearContents()
$ws.Range("G29").Value = 1.4
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = 0
$ws.Range("O29").Value = 0
$ws.Range("P29").Value = 0
$ws.Range("Q29").Value = 0
$ws.Range("R29").Value = ''
$ws.Range("S29").ClearContents()
$ws.Range("T29").ClearContents()
$ws.Range("U29").ClearContents()
$ws.Range("V29").ClearContents()
$ws.Range("W29").ClearContents()
$ws.Range("X29").ClearContents()
$ws.Range("Y29").ClearContents()
$ws.Range("Z29").ClearContents()

# Row 30: 'A 32965-2025' (source row 37)
$ws.Range("A30").Value = 'A 32965-2025'
$ws.Range("B30").Value = 45839
$ws.Range("C30").Value = 46064
$ws.Range("D30").Value = 'STOCKHOLMS LÄN'
$ws.Range("E30").Value = 'SALEM'
$ws.Range("F30").ClearContents()
$ws.Range("G30").Value = 2.8
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 0
$ws.Range("N30").Value = 0
$ws.Range("O30").Value = 0
$ws.Range("P30").Value = 0
$ws.Range("Q30").Value = 0
$ws.Range("R30").Value = ''
$ws.Range("S30").ClearContents()
$ws.Range("T30").ClearContents()
$ws.Range("U30").ClearContents()
$ws.Range("V30").ClearContents()
$ws.Range("W30").ClearContents()
$ws.Range("X30").ClearContents()
$ws.Range("Y30").ClearContents()
$ws.Range("Z30").ClearContents()

# Row 31: 'A 32970-2025' (source row 38)
$ws.Range("A31").Value = 'A 32970-2025'
$ws.Range("B31").Value = 45839
$ws.Range("C31").Value = 46064
$ws.Range("D31").Value = 'STOCKHOLMS LÄN'
$ws.Range("E31").Value = 'SALEM'
$ws.Range("F31").ClearContents()
$ws.Range("G31").Value = 10.5
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 0
$ws.Range("N31").Value = 0
$ws.Range("O31").Value = 0
$ws.Range("P31").Value = 0
$ws.Range("Q31").Value = 0
$ws.Range("R31").Value = ''
$ws.Range("S31").ClearContents()
$ws.Range("T31").ClearContents()
$ws.Range("U31").ClearContents()
$ws.Range("V31").ClearContents()
$ws.Range("W31").ClearContents()
$ws.Range("X31").ClearContents()
$ws.Range("Y31").ClearContents()
$ws.Range("Z31").ClearContents()

# Row 32: 'A 33281-2025' (source row 19)
$ws.Range("A32").Value = 'A 33281-2025'
$ws.Range("B32").Value = 45839
$ws.Range("C32").Value = 46064
$ws.Range("D32").Value = 'STOCKHOLMS LÄN'
$ws.Range("E32").Value = 'SALEM'
$ws.Range("F32").ClearContents()
$ws.Range("G32").Value = 2.6
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = 0
$ws.Range("N32").Value = 0
$ws.Range("O32").Value = 0
$ws.Range("P32").Value = 0
$ws.Range("Q32").Value = 0
$ws.Range("R32").Value = ''
$ws.Range("S32").ClearContents()
$ws.Range("T32").ClearContents()
$ws.Range("U32").ClearContents()
$ws.Range("V32").ClearContents()
$ws.Range("W32").ClearContents()
$ws.Range("X32").ClearContents()
$ws.Range("Y32").ClearContents()
$ws.Range("Z32").ClearContents()

# Row 33: 'A 32953-2025' (source row 43)
$ws.Range("A33").Value = 'A 32953-2025'
$ws.Range("B33").Value = 45839
$ws.Range("C33").Value = 46064
$ws.Range("D33").Value = 'STOCKHOLMS LÄN'
$ws.Range("E33").Value = 'SALEM'
$ws.Range("F33").ClearContents()
$ws.Range("G33").Value = 2.4
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 0
$ws.Range("N33").Value = 0
$ws.Range("O33").Value = 0
$ws.Range("P33").Value = 0
$ws.Range("Q33").Value = 0
$ws.Range("R33").Value = ''
$ws.Range("S33").ClearContents()
$ws.Range("T33").ClearContents()
$ws.Range("U33").ClearContents()
$ws.Range("V33").ClearContents()
$ws.Range("W33").ClearContents()
$ws.Range("X33").ClearContents()
$ws.Range("Y33").ClearContents()
$ws.Range("Z33").ClearContents()

# Row 34: 'A 33470-2025' (source row 32)
$ws.Range("A34").Value = 'A 33470-2025'
$ws.Range("B34").Value = 45839
$ws.Range("C34").Value = 46064
$ws.Range("D34").Value = 'STOCKHOLMS LÄN'
$ws.Range("E34").Value = 'SALEM'
$ws.Range("F34").ClearContents()
$ws.Range("G34").Value = 2.1
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0
$ws.Range("N34").Value = 0
$ws.Range("O34").Value = 0
$ws.Range("P34").Value = 0
$ws.Range("Q34").Value = 0
$ws.Range("R34").Value = ''
$ws.Range("S34").ClearContents()
$ws.Range("T34").ClearContents()
$ws.Range("U34").ClearContents()
$ws.Range("V34").ClearContents()
$ws.Range("W34").ClearContents()
$ws.Range("X34").ClearContents()
$ws.Range("Y34").ClearContents()
$ws.Range("Z34").ClearContents()

# Row 35: 'A 33441-2025' (source row 31)
$ws.Range("A35").Value = 'A 33441-2025'
$ws.Range("B35").Value = 45839
$ws.Range("C35").Value = 46064
$ws.Range("D35").Value = 'STOCKHOLMS LÄN'
$ws.Range("E35").Value = 'SALEM'
$ws.Range("F35").ClearContents()
$ws.Range("G35").Value = 1.7
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = 0
$ws.Range("N35").Value = 0
$ws.Range("O35").Value = 0
$ws.Range("P35").Value = 0
$ws.Range("Q35").Value = 0
$ws.Range("R35").Value = ''
$ws.Range("S35").ClearContents()
$ws.Range("T35").ClearContents()
$ws.Range("U35").ClearContents()
$ws.Range("V35").ClearContents()
$ws.Range("W35").ClearContents()
$ws.Range("X35").ClearContents()
$ws.Range("Y35").ClearContents()
$ws.Range("Z35").ClearContents()

# Row 36: 'A 33428-2025' (source row 20)
$ws.Range("A36").Value = 'A 33428-2025'
$ws.Range("B36").Value = 45839
$ws.Range("C36").Value = 46064
$ws.Range("D36").Value = 'STOCKHOLMS LÄN'
$ws.Range("E36").Value = 'SALEM'
$ws.Range("F36").ClearContents()
$ws.Range("G36").Value = 2.2
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = 0
$ws.Range("N36").Value = 0
$ws.Range("O36").Value = 0
$ws.Range("P36").Value = 0
$ws.Range("Q36").Value = 0
$ws.Range("R36").Value = ''
$ws.Range("S36").ClearContents()
$ws.Range("T36").ClearContents()
$ws.Range("U36").ClearContents()
$ws.Range("V36").ClearContents()
$ws.Range("W36").ClearContents()
$ws.Range("X36").ClearContents()
$ws.Range("Y36").ClearContents()
$ws.Range("Z36").ClearContents()

# Row 37: 'A 33439-2025' (source row 39)
$ws.Range("A37").Value = 'A 33439-2025'
$ws.Range("B37").Value = 45839
$ws.Range("C37").Value = 46064
$ws.Range("D37").Value = 'STOCKHOLMS LÄN'
$ws.Range("E37").Value = 'SALEM'
$ws.Range("F37").ClearContents()
$ws.Range("G37").Value = 2.9
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = 0
$ws.Range("N37").Value = 0
$ws.Range("O37").Value = 0
$ws.Range("P37").Value = 0
$ws.Range("Q37").Value = 0
$ws.Range("R37").Value = ''
$ws.Range("S37").ClearContents()
$ws.Range("T37").ClearContents()
$ws.Range("U37").ClearContents()
$ws.Range("V37").ClearContents()
$ws.Range("W37").ClearContents()
$ws.Range("X37").ClearContents()
$ws.Range("Y37").ClearContents()
$ws.Range("Z37").ClearContents()

# Row 38: 'A 33458-2025' (source row 40)
$ws.Range("A38").Value = 'A 33458-2025'
$ws.Range("B38").Value = 45839
$ws.Range("C38").Value = 46064
$ws.Range("D38").Value = 'STOCKHOLMS LÄN'
$ws.Range("E38").Value = 'SALEM'
$ws.Range("F38").ClearContents()
$ws.Range("G38").Value = 0.8
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 0
$ws.Range("N38").Value = 0
$ws.Range("O38").Value = 0
$ws.Range("P38").Value = 0
$ws.Range("Q38").Value = 0
$ws.Range("R38").Value = ''
$ws.Range("S38").ClearContents()
$ws.Range("T38").ClearContents()
$ws.Range("U38").ClearContents()
$ws.Range("V38").ClearContents()
$ws.Range("W38").ClearContents()
$ws.Range("X38").ClearContents()
$ws.Range("Y38").ClearContents()
$ws.Range("Z38").ClearContents()

# Row 39: 'A 33460-2025' (source row 41)
$ws.Range("A39").Value = 'A 33460-2025'
$ws.Range("B39").Value = 45839
$ws.Range("C39").Value = 46064
$ws.Range("D39").Value = 'STOCKHOLMS LÄN'
$ws.Range("E39").Value = 'SALEM'
$ws.Range("F39").ClearContents()
$ws.Range("G39").Value = 2.1
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = 0
$ws.Range("N39").Value = 0
$ws.Range("O39").Value = 0
$ws.Range("P39").Value = 0
$ws.Range("Q39").Value = 0
$ws.Range("R39").Value = ''
$ws.Range("S39").ClearContents()
$ws.Range("T39").ClearContents()
$ws.Range("U39").ClearContents()
$ws.Range("V39").ClearContents()
$ws.Range("W39").ClearContents()
$ws.Range("X39").ClearContents()
$ws.Range("Y39").ClearContents()
$ws.Range("Z39").ClearContents()

# Row 40: 'A 33465-2025' (source row 42)
$ws.Range("A40").Value = 'A 33465-2025'
$ws.Range("B40").Value = 45839
$ws.Range("C40").Value = 46064
$ws.Range("D40").Value = 'STOCKHOLMS LÄN'
$ws.Range("E40").Value = 'SALEM'
$ws.Range("F40").ClearContents()
$ws.Range("G40").Value = 1.6
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = 0
$ws.Range("N40").Value = 0
$ws.Range("O40").Value = 0
$ws.Range("P40").Value = 0
$ws.Range("Q40").Value = 0
$ws.Range("R40").Value = ''
$ws.Range("S40").ClearContents()
$ws.Range("T40").ClearContents()
$ws.Range("U40").ClearContents()
$ws.Range("V40").ClearContents()
$ws.Range("W40").ClearContents()
$ws.Range("X40").ClearContents()
$ws.Range("Y40").ClearContents()
$ws.Range("Z40").ClearContents()

# Row 41: 'A 32946-2025' (source row 21)
$ws.Range("A41").Value = 'A 32946-2025'
$ws.Range("B41").Value = 45839
$ws.Range("C41").Value = 46064
$ws.Range("D41").Value = 'STOCKHOLMS LÄN'
$ws.Range("E41").Value = 'SALEM'
$ws.Range("F41").ClearContents()
$ws.Range("G41").Value = 1.6
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = 0
$ws.Range("N41").Value = 0
$ws.Range("O41").Value = 0
$ws.Range("P41").Value = 0
$ws.Range("Q41").Value = 0
$ws.Range("R41").Value = ''
$ws.Range("S41").ClearContents()
$ws.Range("T41").ClearContents()
$ws.Range("U41").ClearContents()
$ws.Range("V41").ClearContents()
$ws.Range("W41").ClearContents()
$ws.Range("X41").ClearContents()
$ws.Range("Y41").ClearContents()
$ws.Range("Z41").ClearContents()

# Row 42: 'A 33421-2025' (source row 27)
$ws.Range("A42").Value = 'A 33421-2025'
$ws.Range("B42").Value = 45839
$ws.Range("C42").Value = 46064
$ws.Range("D42").Value = 'STOCKHOLMS LÄN'
$ws.Range("E42").Value = 'SALEM'
$ws.Range("F42").ClearContents()
$ws.Range("G42").Value = 1
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 0
$ws.Range("N42").Value = 0
$ws.Range("O42").Value = 0
$ws.Range("P42").Value = 0
$ws.Range("Q42").Value = 0
$ws.Range("R42").Value = ''
$ws.Range("S42").ClearContents()
$ws.Range("T42").ClearContents()
$ws.Range("U42").ClearContents()
$ws.Range("V42").ClearContents()
$ws.Range("W42").ClearContents()
$ws.Range("X42").ClearContents()
$ws.Range("Y42").ClearContents()
$ws.Range("Z42").ClearContents()

# Row 43: 'A 33446-2025' (source row 26)
$ws.Range("A43").Value = 'A 33446-2025'
$ws.Range("B43").Value = 45839
$ws.Range("C43").Value = 46064
$ws.Range("D43").Value = 'STOCKHOLMS LÄN'
$ws.Range("E43").Value = 'SALEM'
$ws.Range("F43").ClearContents()
$ws.Range("G43").Value = 9.4
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = 0
$ws.Range("N43").Value = 0
$ws.Range("O43").Value = 0
$ws.Range("P43").Value = 0
$ws.Range("Q43").Value = 0
$ws.Range("R43").Value = ''
$ws.Range("S43").ClearContents()
$ws.Range("T43").ClearContents()
$ws.Range("U43").ClearContents()
$ws.Range("V43").ClearContents()
$ws.Range("W43").ClearContents()
$ws.Range("X43").ClearContents()
$ws.Range("Y43").ClearContents()
$ws.Range("Z43").ClearContents()

# Row 44: 'A 33473-2025' (source row 33)
$ws.Range("A44").Value = 'A 33473-2025'
$ws.Range("B44").Value = 45839
$ws.Range("C44").Value = 46064
$ws.Range("D44").Value = 'STOCKHOLMS LÄN'
$ws.Range("E44").Value = 'SALEM'
$ws.Range("F44").ClearContents()
$ws.Range("G44").Value = 5.7
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = 0
$ws.Range("N44").Value = 0
$ws.Range("O44").Value = 0
$ws.Range("P44").Value = 0
$ws.Range("Q44").Value = 0
$ws.Range("R44").Value = ''
$ws.Range("S44").ClearContents()
$ws.Range("T44").ClearContents()
$ws.Range("U44").ClearContents()
$ws.Range("V44").ClearContents()
$ws.Range("W44").ClearContents()
$ws.Range("X44").ClearContents()
$ws.Range("Y44").ClearContents()
$ws.Range("Z44").ClearContents()

# Row 45: 'A 2990-2023' (source row 15)
$ws.Range("A45").Value = 'A 2990-2023'
$ws.Range("B45").Value = 44945
$ws.Range("C45").Value = 46064
$ws.Range("D45").Value = 'STOCKHOLMS LÄN'
$ws.Range("E45").Value = 'SALEM'
$ws.Range("F45").Value = 'Kommuner'
$ws.Range("G45").Value = 4.1
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = 0
$ws.Range("N45").Value = 0
$ws.Range("O45").Value = 0
$ws.Range("P45").Value = 0
$ws.Range("Q45").Value = 0
$ws.Range("R45").Value = ''
$ws.Range("S45").ClearContents()
$ws.Range("T45").ClearContents()
$ws.Range("U45").ClearContents()
$ws.Range("V45").ClearContents()
$ws.Range("W45").ClearContents()
$ws.Range("X45").ClearContents()
$ws.Range("Y45").ClearContents()
$ws.Range("Z45").ClearContents()

# Row 46: 'A 29982-2025' (source row 44)
$ws.Range("A46").Value = 'A 29982-2025'
$ws.Range("B46").Value = 45826
$ws.Range("C46").Value = 46064
$ws.Range("D46").Value = 'STOCKHOLMS LÄN'
$ws.Range("E46").Value = 'SALEM'
$ws.Range("F46").ClearContents()
$ws.Range("G46").Value = 1.5
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = 0
$ws.Range("N46").Value = 0
$ws.Range("O46").Value = 0
$ws.Range("P46").Value = 0
$ws.Range("Q46").Value = 0
$ws.Range("R46").Value = ''
$ws.Range("S46").ClearContents()
$ws.Range("T46").ClearContents()
$ws.Range("U46").ClearContents()
$ws.Range("V46").ClearContents()
$ws.Range("W46").ClearContents()
$ws.Range("X46").ClearContents()
$ws.Range("Y46").ClearContents()
$ws.Range("Z46").ClearContents()

# Row 47: 'A 60731-2025' (source row 50)
$ws.Range("A47").Value = 'A 60731-2025'
$ws.Range("B47").Value = 45995
$ws.Range("C47").Value = 46064
$ws.Range("D47").Value = 'STOCKHOLMS LÄN'
$ws.Range("E47").Value = 'SALEM'
$ws.Range("F47").ClearContents()
$ws.Range("G47").Value = 0.8
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = 0
$ws.Range("N47").Value = 0
$ws.Range("O47").Value = 0
$ws.Range("P47").Value = 0
$ws.Range("Q47").Value = 0
$ws.Range("R47").Value = ''
$ws.Range("S47").ClearContents()
$ws.Range("T47").ClearContents()
$ws.Range("U47").ClearContents()
$ws.Range("V47").ClearContents()
$ws.Range("W47").ClearContents()
$ws.Range("X47").ClearContents()
$ws.Range("Y47").ClearContents()
$ws.Range("Z47").ClearContents()

# Row 48: 'A 60733-2025' (source row 49)
$ws.Range("A48").Value = 'A 60733-2025'
$ws.Range("B48").Value = 45995
$ws.Range("C48").Value = 46064
$ws.Range("D48").Value = 'STOCKHOLMS LÄN'
$ws.Range("E48").Value = 'SALEM'
$ws.Range("F48").ClearContents()
$ws.Range("G48").Value = 1.7
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = 0
$ws.Range("N48").Value = 0
$ws.Range("O48").Value = 0
$ws.Range("P48").Value = 0
$ws.Range("Q48").Value = 0
$ws.Range("R48").Value = ''
$ws.Range("S48").ClearContents()
$ws.Range("T48").ClearContents()
$ws.Range("U48").ClearContents()
$ws.Range("V48").ClearContents()
$ws.Range("W48").ClearContents()
$ws.Range("X48").ClearContents()
$ws.Range("Y48").ClearContents()
$ws.Range("Z48").ClearContents()

# Row 49: 'A 33246-2021' (source row 52)
$ws.Range("A49").Value = 'A 33246-2021'
$ws.Range("B49").Value = 44377
$ws.Range("C49").Value = 46064
$ws.Range("D49").Value = 'STOCKHOLMS LÄN'
$ws.Range("E49").Value = 'SALEM'
$ws.Range("F49").ClearContents()
$ws.Range("G49").Value = 3.4
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = 0
$ws.Range("N49").Value = 0
$ws.Range("O49").Value = 0
$ws.Range("P49").Value = 0
$ws.Range("Q49").Value = 0
$ws.Range("R49").Value = ''
$ws.Range("S49").ClearContents()
$ws.Range("T49").ClearContents()
$ws.Range("U49").ClearContents()
$ws.Range("V49").ClearContents()
$ws.Range("W49").ClearContents()
$ws.Range("X49").ClearContents()
$ws.Range("Y49").ClearContents()
$ws.Range("Z49").ClearContents()

# Row 50: 'A 46919-2023' (source row 47)
$ws.Range("A50").Value = 'A 46919-2023'
$ws.Range("B50").Value = 45201
$ws.Range("C50").Value = 46064
$ws.Range("D50").Value = 'STOCKHOLMS LÄN'
$ws.Range("E50").Value = 'SALEM'
$ws.Range("F50").Value = 'Kommuner'
$ws.Range("G50").Value = 1.3
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = 0
$ws.Range("N50").Value = 0
$ws.Range("O50").Value = 0
$ws.Range("P50").Value = 0
$ws.Range("Q50").Value = 0
$ws.Range("R50").Value = ''
$ws.Range("S50").ClearContents()
$ws.Range("T50").ClearContents()
$ws.Range("U50").ClearContents()
$ws.Range("V50").ClearContents()
$ws.Range("W50").ClearContents()
$ws.Range("X50").ClearContents()
$ws.Range("Y50").ClearContents()
$ws.Range("Z50").ClearContents()

# Row 51: 'A 504-2023' (source row 51)
$ws.Range("A51").Value = 'A 504-2023'
$ws.Range("B51").Value = 44930
$ws.Range("C51").Value = 46064
$ws.Range("D51").Value = 'STOCKHOLMS LÄN'
$ws.Range("E51").Value = 'SALEM'
$ws.Range("F51").Value = 'Kommuner'
$ws.Range("G51").Value = 2
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = 0
$ws.Range("N51").Value = 0
$ws.Range("O51").Value = 0
$ws.Range("P51").Value = 0
$ws.Range("Q51").Value = 0
$ws.Range("R51").Value = ''
$ws.Range("S51").ClearContents()
$ws.Range("T51").ClearContents()
$ws.Range("U51").ClearContents()
$ws.Range("V51").ClearContents()
$ws.Range("W51").ClearContents()
$ws.Range("X51").ClearContents()
$ws.Range("Y51").ClearContents()
$ws.Range("Z51").ClearContents()

# Row 52: 'A 58592-2023' (source row 18)
$ws.Range("A52").Value = 'A 58592-2023'
$ws.Range("B52").Value = 45251
$ws.Range("C52").Value = 46064
$ws.Range("D52").Value = 'STOCKHOLMS LÄN'
$ws.Range("E52").Value = 'SALEM'
$ws.Range("F52").ClearContents()
$ws.Range("G52").Value = 2.1
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = 0
$ws.Range("N52").Value = 0
$ws.Range("O52").Value = 0
$ws.Range("P52").Value = 0
$ws.Range("Q52").Value = 0
$ws.Range("R52").Value = ''
$ws.Range("S52").ClearContents()
$ws.Range("T52").ClearContents()
$ws.Range("U52").ClearContents()
$ws.Range("V52").ClearContents()
$ws.Range("W52").ClearContents()
$ws.Range("X52").ClearContents()
$ws.Range("Y52").ClearContents()
$ws.Range("Z52").ClearContents()

# Row 53: 'A 61336-2022' (source row 46)
$ws.Range("A53").Value = 'A 61336-2022'
$ws.Range("B53").Value = 44915
$ws.Range("C53").Value = 46064
$ws.Range("D53").Value = 'STOCKHOLMS LÄN'
$ws.Range("E53").Value = 'SALEM'
$ws.Range("F53").Value = 'Kommuner'
$ws.Range("G53").Value = 5.6
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = 0
$ws.Range("N53").Value = 0
$ws.Range("O53").Value = 0
$ws.Range("P53").Value = 0
$ws.Range("Q53").Value = 0
$ws.Range("R53").Value = ''
$ws.Range("S53").ClearContents()
$ws.Range("T53").ClearContents()
$ws.Range("U53").ClearContents()
$ws.Range("V53").ClearContents()
$ws.Range("W53").ClearContents()
$ws.Range("X53").ClearContents()
$ws.Range("Y53").ClearContents()
$ws.Range("Z53").ClearContents()

# Row 54: 'A 23992-2023' (source row 55)
$ws.Range("A54").Value = 'A 23992-2023'
$ws.Range("B54").Value = 45078
$ws.Range("C54").Value = 46064
$ws.Range("D54").Value = 'STOCKHOLMS LÄN'
$ws.Range("E54").Value = 'SALEM'
$ws.Range("F54").Value = 'Kommuner'
$ws.Range("G54").Value = 0.7
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = 0
$ws.Range("N54").Value = 0
$ws.Range("O54").Value = 0
$ws.Range("P54").Value = 0
$ws.Range("Q54").Value = 0
$ws.Range("R54").Value = ''
$ws.Range("S54").ClearContents()
$ws.Range("T54").ClearContents()
$ws.Range("U54").ClearContents()
$ws.Range("V54").ClearContents()
$ws.Range("W54").ClearContents()
$ws.Range("X54").ClearContents()
$ws.Range("Y54").ClearContents()
$ws.Range("Z54").ClearContents()

# Row 55: 'A 9533-2024' (source row 54)
$ws.Range("A55").Value = 'A 9533-2024'
$ws.Range("B55").Value = 45359.69105324074
$ws.Range("C55").Value = 46064
$ws.Range("D55").Value = 'STOCKHOLMS LÄN'
$ws.Range("E55").Value = 'SALEM'
$ws.Range("F55").Value = 'Kommuner'
$ws.Range("G55").Value = 2.6
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = 0
$ws.Range("N55").Value = 0
$ws.Range("O55").Value = 0
$ws.Range("P55").Value = 0
$ws.Range("Q55").Value = 0
$ws.Range("R55").Value = ''
$ws.Range("S55").ClearContents()
$ws.Range("T55").ClearContents()
$ws.Range("U55").ClearContents()
$ws.Range("V55").ClearContents()
$ws.Range("W55").ClearContents()
$ws.Range("X55").ClearContents()
$ws.Range("Y55").ClearContents()
$ws.Range("Z55").ClearContents()

# Row 56: 'A 28843-2021' (source row 53)
$ws.Range("A56").Value = 'A 28843-2021'
$ws.Range("B56").Value = 44357
$ws.Range("C56").Value = 46064
$ws.Range("D56").Value = 'STOCKHOLMS LÄN'
$ws.Range("E56").Value = 'SALEM'
$ws.Range("F56").ClearContents()
$ws.Range("G56").Value = 1.3
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = 0
$ws.Range("N56").Value = 0
$ws.Range("O56").Value = 0
$ws.Range("P56").Value = 0
$ws.Range("Q56").Value = 0
$ws.Range("R56").Value = ''
$ws.Range("S56").ClearContents()
$ws.Range("T56").ClearContents()
$ws.Range("U56").ClearContents()
$ws.Range("V56").ClearContents()
$ws.Range("W56").ClearContents()
$ws.Range("X56").ClearContents()
$ws.Range("Y56").ClearContents()
$ws.Range("Z56").ClearContents()

